# Apply the "Office Theme" colour palette to the presentation's active theme
# (ppt/theme/theme1.xml), replacing the current "Integral" / "Red Violet"
# palette -- equivalent to picking a different built-in Design/Theme in the
# PowerPoint Design gallery.

function MyRGB($r, $g, $b) {
    $v = $r + ($g * 256) + ($b * 65536)
    return $v
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# ThemeColorScheme item order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink
$cs.Item(1).RGB  = (MyRGB 0x00 0x00 0x00)   # dk1
$cs.Item(2).RGB  = (MyRGB 0xFF 0xFF 0xFF)   # lt1
$cs.Item(3).RGB  = (MyRGB 0x44 0x54 0x6A)   # dk2
$cs.Item(4).RGB  = (MyRGB 0xE7 0xE6 0xE6)   # lt2
$cs.Item(5).RGB  = (MyRGB 0x5B 0x9B 0xD5)   # accent1
$cs.Item(6).RGB  = (MyRGB 0xED 0x7D 0x31)   # accent2
$cs.Item(7).RGB  = (MyRGB 0xA5 0xA5 0xA5)   # accent3
$cs.Item(8).RGB  = (MyRGB 0xFF 0xC0 0x00)   # accent4
$cs.Item(9).RGB  = (MyRGB 0x44 0x72 0xC4)   # accent5
$cs.Item(10).RGB = (MyRGB 0x70 0xAD 0x47)   # accent6
$cs.Item(11).RGB = (MyRGB 0x05 0x63 0xC1)   # hlink
$cs.Item(12).RGB = (MyRGB 0x95 0x4F 0x72)   # folHlink
